$d = $word.ActiveDocument

# 1. Update the heading text and drop the "Bradicardia" paragraph that
#    used to follow it. We replace the heading run's text first, then
#    remove the now-separate "Bradicardia – câmara dupla." paragraph
#    (including its paragraph mark) entirely.
$d.Content.Find.Execute("Implante de Marcapasso Convencional (Astra" + [char]0x2122 + ")", $true, $false, $false, $false, $false, $true, 1, $false, "Marcapasso Dupla Câmara (Astra)", 2)

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Bradicardia " + [char]0x2013 + " câmara dupla.") {
        $p.Range.Delete()
        break
    }
}

# 2. Material list items: drop "™" where present, normalise the
#    separators, and prefix each line with a bullet glyph.
$d.Content.Find.Execute("Gerador " + [char]0x2013 + " Astra" + [char]0x2122, $true, $false, $false, $false, $false, $true, 1, $false, [char]0x2022 + " Gerador Astra", 2)
$d.Content.Find.Execute("Eletrodo Ventricular " + [char]0x2013 + " 5076-52", $true, $false, $false, $false, $false, $true, 1, $false, [char]0x2022 + " Eletrodo 5076-52", 2)
$d.Content.Find.Execute("Eletrodo Atrial " + [char]0x2013 + " 5076-58", $true, $false, $false, $false, $false, $true, 1, $false, [char]0x2022 + " Eletrodo 5076-58", 2)
$d.Content.Find.Execute("Introdutor " + [char]0x2013 + " 2", $true, $false, $false, $false, $false, $true, 1, $false, [char]0x2022 + " Introdutor " + [char]0x2013 + " 2", 2)
